$d = $word.ActiveDocument
$ErrorActionPreference = "Stop"

function Replace-ParagraphByAnchor([string]$Anchor, [string]$Xml) {
    $rng = $d.Content
    $found = $rng.Find.Execute($Anchor, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Anchor not found: $Anchor"
    }
    $p = $rng.Paragraphs(1)
    $prng = $p.Range
    $prng.InsertXML($Xml)
}

$xml1 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w14:paraId="5058DC42" w14:textId="77777777" w:rsidR="00F84099" w:rsidRDefault="00F84099" w:rsidP="00F84099"><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:ind w:left="1080"/></w:pPr><w:r><w:t xml:space="preserve">Decide quindi di premere il tasto “Candidati ora”. Gli viene richiesto di fare login, </w:t></w:r><w:r><w:t xml:space="preserve">tuttavia Giuseppe non possiede un account, perciò clicca su “Registrati come inoccupato”. </w:t></w:r><w:r><w:t xml:space="preserve">A questo punto Giuseppe visualizza </w:t></w:r><w:r><w:t>un modulo da compilare con nome: “Giuseppe”, cognome: “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Rainone</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”, città: “Salerno”, indirizzo: “</w:t></w:r><w:r><w:t>Viale Europa, 44</w:t></w:r><w:r><w:t>”, data di nascita: “</w:t></w:r><w:r><w:t>28/09/1991</w:t></w:r><w:r><w:t>”</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>username:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>“Giuseppe33</w:t></w:r><w:r><w:t>”,</w:t></w:r><w:r><w:t xml:space="preserve"> password: “051Progetto_”</w:t></w:r><w:r><w:t xml:space="preserve"> e carica il proprio curricul</w:t></w:r><w:r><w:t>u</w:t></w:r><w:r><w:t>m: “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Rainone</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-Giuseppe-CV”; spunta la casella per il trattamento dei dati e clicca su “Conferma”. Viene quindi reindirizzato ad una pagina intermedia</w:t></w:r><w:r><w:t xml:space="preserve"> per notificare l’avvenuta registrazione dove gli viene richiesta una conferma per continuare la candidatura. Giuseppe clicca su “Continua” e viene reindirizzato ad un’altra pagina intermedia</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
Replace-ParagraphByAnchor 'Decide quindi di premere il tasto' $xml1

$xml2 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w14:paraId="3D02CD4B" w14:textId="77777777" w:rsidR="00F84099" w:rsidRDefault="00F84099" w:rsidP="00F84099"><w:pPr><w:ind w:left="1080"/></w:pPr><w:r><w:t>Veronica, vedendo che Giuseppe è di Salerno, decide di visualizzare il suo curriculum e clicca su “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Rainone</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">-Giuseppe-CV”, si apre una nuova pagina esterna che le permette di visualizzare il file PDF del candidato. Leggendo il curriculum, decide che Giuseppe sarebbe un buon candidato, in quanto soddisfa i requisiti per l’occupazione. A questo punto Veronica torna alla pagina delle </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>candidature e clicca su “Contatta”, le viene quindi mostrato un modulo di compilazione in cui inserisce un titolo</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>: ”Appuntamento</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> per un colloquio – </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Modis</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> S.p.A.” e un corpo: </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
Replace-ParagraphByAnchor 'Veronica, vedendo che Giuseppe' $xml2

$xml3 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w14:paraId="2F8A7F74" w14:textId="77777777" w:rsidR="00F84099" w:rsidRDefault="00F84099" w:rsidP="00F84099"><w:pPr><w:ind w:left="1080"/></w:pPr><w:r><w:t xml:space="preserve">“Gentile Giuseppe </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Rainone</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>,</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">dopo aver analizzato le sue competenze, la riteniamo un candidato adeguato </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>per il</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> ruolo indicato dal nostro annuncio. La invitiamo per un colloquio in data 19/12/2019 presso la nostra sede centrale in Via Torquato Tasso, 34 situata a Salerno.</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">Per ulteriori informazioni potrà contattarci al nostro indirizzo e-mail: </w:t></w:r><w:hyperlink r:id="rId5" w:history="1"><w:r w:rsidRPr="00796733"><w:rPr><w:rStyle w:val="Collegamentoipertestuale"/></w:rPr><w:t>modis.salerno@gmail.com</w:t></w:r></w:hyperlink><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
Replace-ParagraphByAnchor '“Gentile Giuseppe' $xml3

$xml4 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w14:paraId="728B8AB3" w14:textId="30223B26" w:rsidR="00F84099" w:rsidRDefault="00F84099" w:rsidP="00F84099"><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:ind w:left="1080"/></w:pPr><w:r><w:t>L’azienda “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>BusItalia</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">” ha percepito un bonus dalla Regione Campania ed ha ampliato i suoi servizi di trasporto in termini di nuove corse per tutta la Regione. Ciò ha quindi generato l’esigenza di assumere nuovi autisti. Francesco, il responsabile delle human </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>resource</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> di “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>BusItalia</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”, viene incaricato di inserire un annuncio sulla nostra piattaforma. Francesco quindi si collega al sito ed effettua la registrazione sotto forma di profilo aziendale cliccando su “Registrati come azienda”. Gli compare una nuova pagina dove inserisce nome dell’azienda: “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>BusItalia</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> S.p.A.”, immagine del logo dell’</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>azienda</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>:”bu</w:t></w:r><w:r w:rsidR="002A33F5"><w:t>s</w:t></w:r><w:r><w:t>italia_logo.jpg</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>”, username: “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>BusItaliaRoma</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">”, indirizzo: “Via V Maggio, 123, Roma”, data di fondazione: “19 Maggio 2011”, descrizione sintetica: “Azienda di trasporto pubblico locale.”, numero di dipendenti attuali: “3639”, e-mail: “busitaliaroma@gmail.com”, password: “Informazioni_01”. Una volta compilata correttamente la scheda, conferma il modulo e viene reindirizzato alla pagina di avvenuta registrazione. A questo punto Francesco, accedendo alla propria pagina, clicca su “Nuova </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">inserzione” e gli viene mostrata la pagina dedicata alla creazione dell’annuncio. Compila il modulo inserendo un titolo: “Cercasi autisti per autobus”, tag: “Trasporto”, descrizione: “Azienda leader dei trasporti in Italia, cerca nuovi profili da inserire nel proprio organico”, requisiti: “Richiesto il possesso della patente CQC da almeno 5 anni con esperienza </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>regressa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> nel settore di almeno 2 anni.”, contratto: “Full-time”. Francesco revisiona la scheda, clicca “Pubblica” e viene reindirizzato alla pagina di avvenuta pubblicazione.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
Replace-ParagraphByAnchor '’azienda “' $xml4

$xml5 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w14:paraId="6E3E354B" w14:textId="77777777" w:rsidR="00F84099" w:rsidRDefault="00F84099" w:rsidP="00F84099"><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Security</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
Replace-ParagraphByAnchor 'Security' $xml5
